# Tasks.xlsx edit script
# Commit: "Updated with Graphs and Examples for augmentation and Language
#          translation. Tweaked augmentation also. Might need to rerun it
#          and maybe it will improve accuracy"
#
# Changes applied:
#   1. Status column updates (D8, D16, D21): mark three tasks "Completed"
#      (previously "Almost Complete" / "Pending" / "Pending").
#   2. Remove the AutoFilter from Table1 (no longer needed/shown).
#   3. Update the active selection on Sheet1 from B24 to D14 (reflects
#      where the author was last working/reviewing).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update task statuses to "Completed"
$ws.Range("D8").Value  = "Completed"
$ws.Range("D16").Value = "Completed"
$ws.Range("D21").Value = "Completed"

# 2. Turn off the table's AutoFilter dropdowns / remove the autoFilter element
$table = $ws.ListObjects.Item("Table1")
$table.ShowAutoFilter = $false

# 3. Move the active selection to D14
$ws.Range("D14").Select()
